# peserta.xlsx template enhancement:
# - split "No Rekening" duplication into NIK / NPWP columns (mirrors NIP)
# - rename "No Rekening" -> "Nomor Rekening"
# - fix the Golongan/Ruang sample values
# - move the "Referensi Golongan/Ruang" helper list from column L to column O
#   and repoint the data validation list at it
# - add a note on the Golongan/Ruang header listing the allowed values
# - freeze the first 3 (ID/Nama/NIP) columns + header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. make room: two new columns right after NIP (C) for NIK / NPWP ---
$ws.Columns("D:E").Insert()

# --- 2. make room: one new column so the reference list ends up in O (was L) ---
$ws.Columns("N:N").Insert()

# --- 3. header row ---
$ws.Range("D1").Value = "NIK"
$ws.Range("E1").Value = "NPWP"
$ws.Range("K1").Value = "Nomor Rekening"

# --- 4. body rows: NIK / NPWP mirror the NIP column, like the template's other id cols ---
$ws.Range("D2").Value = $ws.Range("C2").Text
$ws.Range("D3").Value = $ws.Range("C3").Text
$ws.Range("D4").Value = $ws.Range("C4").Text
$ws.Range("E2").Value = $ws.Range("C2").Text
$ws.Range("E3").Value = $ws.Range("C3").Text
$ws.Range("E4").Value = $ws.Range("C4").Text

# copy NIP's number-as-text formatting (quote-prefixed, bordered) onto the new cells
$ws.Range("C2").Copy()
$ws.Range("D2:E2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("D3:E3").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("D4:E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. fix the sample Golongan/Ruang values (column F after the insert) ---
$ws.Range("F2").Value = "I/D"
$ws.Range("F3").Value = "III/D"
$ws.Range("F4").Value = "III/c"

# --- 6. strip the old bordered-list look off the relocated reference list (now col O) ---
$ws.Range("O1:O18").ClearFormats()
$ws.Range("O1").Value = "Referensi Golongan/Ruang"

# --- 7. data validation: drop the stale rule, add the one pointing at O2:O18 ---
$ws.Cells.Validation.Delete()
$rngF = $ws.Range("F1:F1048576")
$rngF.Validation.Add(3, 1, 1, "=`$O`$2:`$O`$18")
$rngF.Validation.ErrorTitle = "Golongan/Ruang"
$rngF.Validation.ErrorMessage = "Data tidak benar"
$rngF.Validation.IgnoreBlank = $true
$rngF.Validation.InCellDropdown = $true
$rngF.Validation.ShowInput = $true
$rngF.Validation.ShowError = $true

# --- 8. note listing the allowed values, on the Golongan/Ruang header ---
$ws.Range("F1").AddComment("I/A, I/B, I/C, I/D, II/A, II/B, II/C, II/D, III/A, III/B, III/C, III/D, IV/A, IV/B, IV/C, IV/D, IV/E")

# --- 9. column widths to fit the new headers/content ---
$ws.Columns("F:F").ColumnWidth = 15.666666666666666
$ws.Columns("G:G").ColumnWidth = 13.666666666666666
$ws.Columns("I:I").ColumnWidth = 13.330729166666666
$ws.Columns("K:K").ColumnWidth = 10.666666666666666
$ws.Columns("O:O").ColumnWidth = 22.830729166666668

# --- 10. freeze the ID/Nama/NIP columns + header row ---
$ws.Range("D2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E15").Select()

Write-Host "done"
